# Rename the "FirstFinalScore" / "SecondFinalScore" column headers to a
# unified "FinalScore" header on both department sheets.

$wb = $excel.ActiveWorkbook

$wsInnovation = $wb.Worksheets.Item("Innovation Department")
$wsInnovation.Range("H1").Value = "FinalScore"

$wsProcurement = $wb.Worksheets.Item("Procurement Department")
$wsProcurement.Range("K1").Value = "FinalScore"
